$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45: the C45 "NA" placeholder becomes blank now that a real page number exists
$ws.Range("C45").Value = ""

# Append the new rows produced by the latest script run
$ws.Range("A46").Value = "'2025-04-11"
$ws.Range("A46").Style = "Normal"
$ws.Range("B46").Value = "organismes de quarantaine"
$ws.Range("C46").Value = 7
$ws.Range("D46").Value = 1

$ws.Range("A47").Value = "'2025-04-11"
$ws.Range("A47").Style = "Normal"
$ws.Range("B47").Value = "organismes de quarantaine"
$ws.Range("C47").Value = 112
$ws.Range("D47").Value = 2

$ws.Range("A48").Value = "'2025-04-11"
$ws.Range("A48").Style = "Normal"
$ws.Range("B48").Value = "développement durable"
$ws.Range("C48").Value = 162
$ws.Range("D48").Value = 1
